$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value as plain TEXT (not auto-converted to a number),
# preserving exact digit formatting (leading/trailing zeros, etc.), and
# avoid leaving an extra "quote prefix" style behind on the cell.
function Set-TextValue($cellRef, $val) {
    $ws.Range($cellRef).Value = "'" + $val
    $ws.Range($cellRef).Style = "Normal"
}

# --- Price (column D) updates ---
Set-TextValue "D2"  "247.12"
Set-TextValue "D3"  "22.04"
Set-TextValue "D4"  "5.455"
Set-TextValue "D6"  "3.400"
Set-TextValue "D7"  "6.342"
Set-TextValue "D8"  "0.8186"
Set-TextValue "D9"  "0.9790"
Set-TextValue "D13" "0.02997"
Set-TextValue "D14" "4.167"
Set-TextValue "D15" "0.09422"
Set-TextValue "D16" "0.001596"
Set-TextValue "D17" "0.04823"
Set-TextValue "D18" "0.0005848"
Set-TextValue "D19" "0.006184"
Set-TextValue "D20" "0.004132"
Set-TextValue "D21" "0.0009985"
Set-TextValue "D22" "0.0001500"
Set-TextValue "D23" "3.762"
Set-TextValue "D24" "2.218"
Set-TextValue "D25" "0.3257"
Set-TextValue "D27" "0.0003998"
Set-TextValue "D40" "0.03893"

# --- Row 9: also update the Volume(1h) label text ---
Set-TextValue "E9" "8FTXTokenFTTBestin24h"

# --- Rows 41-43: the coin list got re-ranked (Kick/BKEX/CEJI rotated) ---
# Row 41 : BKEXToken -> KickToken
Set-TextValue "B41" "KickToken"
Set-TextValue "C41" "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
Set-TextValue "D41" "0.006460"
Set-TextValue "E41" "40KickTokenKICK"

# Row 42 : CEJI -> BKEXToken
Set-TextValue "B42" "BKEXToken"
Set-TextValue "C42" "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
Set-TextValue "D42" "0.1076"
Set-TextValue "E42" "41BKEXTokenBKK"

# Row 43 : KickToken -> CEJI
Set-TextValue "B43" "CEJI"
Set-TextValue "C43" "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
Set-TextValue "D43" "0.003000"
Set-TextValue "E43" "42CEJICEJI"

# --- Remaining price (column D) updates ---
Set-TextValue "D44" "0.006497"
Set-TextValue "D45" "0.00005591"
Set-TextValue "D47" "0.3799"
